$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting existing rows 59:131 down to 60:132
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new record's data
$ws.Range("A59").Value = 4
$ws.Range("B59").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C59").Value = "Los Lagos"
$ws.Range("D59").Value = 44483
$ws.Range("E59").Value = 10
$ws.Range("F59").Value = 100112032
$ws.Range("G59").Value = "Zapallo italiano"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 200
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 15000
$ws.Range("M59").Value = 15000
$ws.Range("N59").Value = "`$/caja 50 unidades"
$ws.Range("O59").Value = "Región de Arica y Parinacota"
$ws.Range("P59").Value = 300
$ws.Range("Q59").Value = 50
$ws.Range("R59").Value = "Hortaliza"
